$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "26.943.45"
$ws.Range("E2").Value = "  +1.58%  "

$ws.Range("D3").Value = "1.847.27"
$ws.Range("E3").Value = "  +1.84%  "

Set-TextValue $ws.Range("D4") "1.008"
$ws.Range("E4").Value = "  +0.49%  "

Set-TextValue $ws.Range("D5") "309.66"
$ws.Range("E5").Value = "  +1.17%  "

Set-TextValue $ws.Range("D6") "1.007"
$ws.Range("E6").Value = "  +0.39%  "

Set-TextValue $ws.Range("D7") "0.4683"
$ws.Range("E7").Value = "  +3.31%  "

Set-TextValue $ws.Range("D8") "0.3662"
$ws.Range("E8").Value = "  +2.01%  "

Set-TextValue $ws.Range("D9") "0.07158"
$ws.Range("E9").Value = "  +0.84%  "

Set-TextValue $ws.Range("D10") "0.9288"
$ws.Range("E10").Value = "  +3.58%  "

Set-TextValue $ws.Range("D11") "19.60"
$ws.Range("E11").Value = "  +1.23%  "

Set-TextValue $ws.Range("D12") "0.07704"
$ws.Range("E12").Value = "  -0.94%  "

$ws.Range("D13").Value = "1.843.55"
$ws.Range("E13").Value = "  +0.93%  "

Set-TextValue $ws.Range("D14") "5.290"
$ws.Range("E14").Value = "  +0.31%  "

Set-TextValue $ws.Range("D15") "6.410"
$ws.Range("E15").Value = "  +1.29%  "

Set-TextValue $ws.Range("D16") "88.43"
$ws.Range("E16").Value = "  +3.76%  "

Set-TextValue $ws.Range("D17") "1.009"
$ws.Range("E17").Value = "  +0.44%  "

Set-TextValue $ws.Range("D18") "0.000008626"
$ws.Range("E18").Value = "  +0.31%  "

$ws.Range("E19").Value = "  +0.38%  "

$ws.Range("D20").Value = "26.979.22"
$ws.Range("E20").Value = "  +1.57%  "

Set-TextValue $ws.Range("D22") "5.027"
$ws.Range("E22").Value = "  +1.18%  "

$ws.Range("E23").Value = "  +0.79%  "

$ws.Range("E24").Value = "  -2.04%  "

Set-TextValue $ws.Range("D25") "152.36"
$ws.Range("E25").Value = "  +0.99%  "

$ws.Range("E26").Value = "  +2.54%  "

Set-TextValue $ws.Range("D27") "2.021"
$ws.Range("E27").Value = "  -1.54%  "

Set-TextValue $ws.Range("D28") "114.52"
$ws.Range("E28").Value = "  +1.88%  "

Set-TextValue $ws.Range("D29") "4.880"
$ws.Range("E29").Value = "  +0.48%  "

Set-TextValue $ws.Range("D30") "0.08859"
$ws.Range("E30").Value = "  +1.59%  "

Set-TextValue $ws.Range("D31") "3.215"
$ws.Range("E31").Value = "  +2.97%  "

Set-TextValue $ws.Range("D32") "1.176"
$ws.Range("E32").Value = "  +5.75%  "

Set-TextValue $ws.Range("D33") "0.7481"
$ws.Range("E33").Value = "  -1.03%  "

Set-TextValue $ws.Range("D34") "2.802"
$ws.Range("E34").Value = "  +2.63%  "

Set-TextValue $ws.Range("D35") "4.477"
$ws.Range("E35").Value = "  +1.01%  "

Set-TextValue $ws.Range("D36") "1.087"
$ws.Range("E36").Value = "  +1.37%  "

Set-TextValue $ws.Range("D38") "2.970"
$ws.Range("E38").Value = "  +2.12%  "

Set-TextValue $ws.Range("D39") "0.05196"
$ws.Range("E39").Value = "  +1.72%  "

Set-TextValue $ws.Range("D40") "0.5208"
$ws.Range("E40").Value = "  +2.17%  "

Set-TextValue $ws.Range("D41") "6.914"
$ws.Range("E41").Value = "  +2.12%  "

Set-TextValue $ws.Range("D42") "0.1516"
$ws.Range("E42").Value = "  +0.49%  "

$ws.Range("E43").Value = "  +1.19%  "

Set-TextValue $ws.Range("D44") "10.56"
$ws.Range("E44").Value = "  +5.67%  "

Set-TextValue $ws.Range("D45") "0.4696"
$ws.Range("E45").Value = "  -0.41%  "

Set-TextValue $ws.Range("D47") "100.24"
$ws.Range("E47").Value = "  -0.93%  "

$ws.Range("E48").Value = "  +1.86%  "

Set-TextValue $ws.Range("D49") "65.84"
$ws.Range("E49").Value = "  +2.97%  "

Set-TextValue $ws.Range("D50") "0.06046"
$ws.Range("E50").Value = "  +1.04%  "

Set-TextValue $ws.Range("D51") "0.8933"
$ws.Range("E51").Value = "  +5.29%  "
